$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert a new worksheet "Planilha2" right after "Planilha1"
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "Planilha2"

# Move (cut) the data rows for products 2-5 (A3:F6) from Planilha1
# into the new Planilha2 sheet, landing at A1:F4
$src = $ws1.Range("A3:F6")
$src.Cut($ws2.Range("A1:F4"))

# The vacated rows in Planilha1 keep getting reformatted: no border,
# underlined font (new style)
$empty = $ws1.Range("A3:F6")
$empty.Borders.LineStyle = -4142
$empty.Font.Underline = $true

# Give Planilha1 an explicit page setup (A4, portrait)
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# Keep the same page margins on the new sheet as on Planilha1
$ws2.PageSetup.LeftMargin = $ws1.PageSetup.LeftMargin
$ws2.PageSetup.RightMargin = $ws1.PageSetup.RightMargin
$ws2.PageSetup.TopMargin = $ws1.PageSetup.TopMargin
$ws2.PageSetup.BottomMargin = $ws1.PageSetup.BottomMargin
$ws2.PageSetup.HeaderMargin = $ws1.PageSetup.HeaderMargin
$ws2.PageSetup.FooterMargin = $ws1.PageSetup.FooterMargin

# Leave the selection on the new sheet roughly where the pasted data landed
$ws2.Range("A1:F4").Select() | Out-Null

# Re-activate Planilha1 as the visible/selected sheet, cursor at B10
$ws1.Activate() | Out-Null
$ws1.Range("B10").Select() | Out-Null
